$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the merged header "Fully bookedDate obtained" into two separate headers
$ws.Range("H1").Value = "Fully booked"
$ws.Range("I1").Value = "Date obtained"
